$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Update the two summary cells (reporting period + download timestamp)
# ------------------------------------------------------------------
$ws.Range("B4").Value = "2024-12-01 ~ 2024-12-31"
$ws.Range("B5").Value = "2025년 01월 10일 01시 32분 28초"

# ------------------------------------------------------------------
# Snapshot the existing cell formatting into scratch cells (columns H-M)
# before we touch any values, so we can restore the styles afterwards.
#   H8:M9   -> two-row alternating style pattern used by data rows (s=7/8 and s=7/9)
#   H38:M38 -> style pattern used by the trailing blank row (s=4)
# ------------------------------------------------------------------
$ws.Range("A8:F9").Copy()
$ws.Range("H8:M9").PasteSpecial(-4122)
$ws.Range("A38:F38").Copy()
$ws.Range("H38:M38").PasteSpecial(-4122)

# Force the whole data range (including the new row 39) to be treated as plain
# text so that date-like strings ("2024-12-31", ...) are not auto-converted to
# Excel date serial numbers when we assign them below.
$ws.Range("A8:F39").NumberFormat = "@"

# ------------------------------------------------------------------
# Write the new values. Rows 8-37 shift the existing Nov data down to Dec,
# row 38 becomes a brand-new data row (2024-12-01), and a new trailing blank
# row is created at row 39.
# ------------------------------------------------------------------
$ws.Range("A8").Value = "2024-12-31"
$ws.Range("B8").Value = "화"
$ws.Range("C8").Value = "50"
$ws.Range("D8").Value = "0"
$ws.Range("E8").Value = "1"
$ws.Range("F8").Value = "49"
$ws.Range("A9").Value = "2024-12-30"
$ws.Range("B9").Value = "월"
$ws.Range("C9").Value = "44"
$ws.Range("D9").Value = "1"
$ws.Range("E9").Value = "0"
$ws.Range("F9").Value = "43"
$ws.Range("A10").Value = "2024-12-29"
$ws.Range("B10").Value = "일"
$ws.Range("C10").Value = "47"
$ws.Range("D10").Value = "1"
$ws.Range("E10").Value = "0"
$ws.Range("F10").Value = "46"
$ws.Range("A11").Value = "2024-12-28"
$ws.Range("B11").Value = "토"
$ws.Range("C11").Value = "73"
$ws.Range("D11").Value = "1"
$ws.Range("E11").Value = "0"
$ws.Range("F11").Value = "72"
$ws.Range("A12").Value = "2024-12-27"
$ws.Range("B12").Value = "금"
$ws.Range("C12").Value = "45"
$ws.Range("D12").Value = "0"
$ws.Range("E12").Value = "0"
$ws.Range("F12").Value = "45"
$ws.Range("A13").Value = "2024-12-26"
$ws.Range("B13").Value = "목"
$ws.Range("C13").Value = "60"
$ws.Range("D13").Value = "8"
$ws.Range("E13").Value = "0"
$ws.Range("F13").Value = "52"
$ws.Range("A14").Value = "2024-12-25"
$ws.Range("B14").Value = "수"
$ws.Range("C14").Value = "47"
$ws.Range("D14").Value = "6"
$ws.Range("E14").Value = "2"
$ws.Range("F14").Value = "39"
$ws.Range("A15").Value = "2024-12-24"
$ws.Range("B15").Value = "화"
$ws.Range("C15").Value = "82"
$ws.Range("D15").Value = "9"
$ws.Range("E15").Value = "2"
$ws.Range("F15").Value = "71"
$ws.Range("A16").Value = "2024-12-23"
$ws.Range("B16").Value = "월"
$ws.Range("C16").Value = "9"
$ws.Range("D16").Value = "0"
$ws.Range("E16").Value = "0"
$ws.Range("F16").Value = "9"
$ws.Range("A17").Value = "2024-12-22"
$ws.Range("B17").Value = "일"
$ws.Range("C17").Value = "3"
$ws.Range("D17").Value = "0"
$ws.Range("E17").Value = "0"
$ws.Range("F17").Value = "3"
$ws.Range("A18").Value = "2024-12-21"
$ws.Range("B18").Value = "토"
$ws.Range("C18").Value = "6"
$ws.Range("D18").Value = "0"
$ws.Range("E18").Value = "0"
$ws.Range("F18").Value = "6"
$ws.Range("A19").Value = "2024-12-20"
$ws.Range("B19").Value = "금"
$ws.Range("C19").Value = "6"
$ws.Range("D19").Value = "0"
$ws.Range("E19").Value = "0"
$ws.Range("F19").Value = "6"
$ws.Range("A20").Value = "2024-12-19"
$ws.Range("B20").Value = "목"
$ws.Range("C20").Value = "8"
$ws.Range("D20").Value = "0"
$ws.Range("E20").Value = "0"
$ws.Range("F20").Value = "8"
$ws.Range("A21").Value = "2024-12-18"
$ws.Range("B21").Value = "수"
$ws.Range("C21").Value = "5"
$ws.Range("D21").Value = "0"
$ws.Range("E21").Value = "0"
$ws.Range("F21").Value = "5"
$ws.Range("A22").Value = "2024-12-17"
$ws.Range("B22").Value = "화"
$ws.Range("C22").Value = "5"
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "0"
$ws.Range("F22").Value = "5"
$ws.Range("A23").Value = "2024-12-16"
$ws.Range("B23").Value = "월"
$ws.Range("C23").Value = "11"
$ws.Range("D23").Value = "0"
$ws.Range("E23").Value = "0"
$ws.Range("F23").Value = "11"
$ws.Range("A24").Value = "2024-12-15"
$ws.Range("B24").Value = "일"
$ws.Range("C24").Value = "3"
$ws.Range("D24").Value = "0"
$ws.Range("E24").Value = "0"
$ws.Range("F24").Value = "3"
$ws.Range("A25").Value = "2024-12-14"
$ws.Range("B25").Value = "토"
$ws.Range("C25").Value = "2"
$ws.Range("D25").Value = "0"
$ws.Range("E25").Value = "0"
$ws.Range("F25").Value = "2"
$ws.Range("A26").Value = "2024-12-13"
$ws.Range("B26").Value = "금"
$ws.Range("C26").Value = "7"
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "0"
$ws.Range("F26").Value = "7"
$ws.Range("A27").Value = "2024-12-12"
$ws.Range("B27").Value = "목"
$ws.Range("C27").Value = "2"
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "0"
$ws.Range("F27").Value = "2"
$ws.Range("A28").Value = "2024-12-11"
$ws.Range("B28").Value = "수"
$ws.Range("C28").Value = "5"
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "0"
$ws.Range("F28").Value = "5"
$ws.Range("A29").Value = "2024-12-10"
$ws.Range("B29").Value = "화"
$ws.Range("C29").Value = "6"
$ws.Range("D29").Value = "0"
$ws.Range("E29").Value = "0"
$ws.Range("F29").Value = "6"
$ws.Range("A30").Value = "2024-12-09"
$ws.Range("B30").Value = "월"
$ws.Range("C30").Value = "13"
$ws.Range("D30").Value = "1"
$ws.Range("E30").Value = "0"
$ws.Range("F30").Value = "12"
$ws.Range("A31").Value = "2024-12-08"
$ws.Range("B31").Value = "일"
$ws.Range("C31").Value = "7"
$ws.Range("D31").Value = "0"
$ws.Range("E31").Value = "0"
$ws.Range("F31").Value = "7"
$ws.Range("A32").Value = "2024-12-07"
$ws.Range("B32").Value = "토"
$ws.Range("C32").Value = "8"
$ws.Range("D32").Value = "1"
$ws.Range("E32").Value = "0"
$ws.Range("F32").Value = "7"
$ws.Range("A33").Value = "2024-12-06"
$ws.Range("B33").Value = "금"
$ws.Range("C33").Value = "12"
$ws.Range("D33").Value = "0"
$ws.Range("E33").Value = "0"
$ws.Range("F33").Value = "12"
$ws.Range("A34").Value = "2024-12-05"
$ws.Range("B34").Value = "목"
$ws.Range("C34").Value = "8"
$ws.Range("D34").Value = "0"
$ws.Range("E34").Value = "0"
$ws.Range("F34").Value = "8"
$ws.Range("A35").Value = "2024-12-04"
$ws.Range("B35").Value = "수"
$ws.Range("C35").Value = "7"
$ws.Range("D35").Value = "2"
$ws.Range("E35").Value = "0"
$ws.Range("F35").Value = "5"
$ws.Range("A36").Value = "2024-12-03"
$ws.Range("B36").Value = "화"
$ws.Range("C36").Value = "10"
$ws.Range("D36").Value = "3"
$ws.Range("E36").Value = "1"
$ws.Range("F36").Value = "6"
$ws.Range("A37").Value = "2024-12-02"
$ws.Range("B37").Value = "월"
$ws.Range("C37").Value = "10"
$ws.Range("D37").Value = "2"
$ws.Range("E37").Value = "0"
$ws.Range("F37").Value = "8"
$ws.Range("A38").Value = "2024-12-01"
$ws.Range("B38").Value = "일"
$ws.Range("C38").Value = "7"
$ws.Range("D38").Value = "0"
$ws.Range("E38").Value = "0"
$ws.Range("F38").Value = "7"

# ------------------------------------------------------------------
# Restore the correct cell styles now that values are populated:
#  - rows 8-37 get the tiled alternating style pattern
#  - row 38 gets the single (even) data-row style
#  - the new row 39 gets the trailing blank-row style
# ------------------------------------------------------------------
$ws.Range("H8:M9").Copy()
$ws.Range("A8:F37").PasteSpecial(-4122)

$ws.Range("H8:M8").Copy()
$ws.Range("A38:F38").PasteSpecial(-4122)

$ws.Range("H38:M38").Copy()
$ws.Range("A39:F39").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Clean up the scratch template cells
# ------------------------------------------------------------------
$ws.Range("H8:M9").Clear()
$ws.Range("H38:M38").Clear()

Write-Output "Edit complete"
